$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: duplicate of row 3 but with Focal_Species (E) = "Badger"
$ws.Range("A3:J3").Copy($ws.Range("A5"))
$ws.Range("E5").Value = "Badger"
$ws.Range("H5").FormulaR1C1 = '=""'

# Row 6: duplicate of row 3 but with Focal_Species (E) = "badger"
$ws.Range("A3:J3").Copy($ws.Range("A6"))
$ws.Range("E6").Value = "badger"
$ws.Range("H6").FormulaR1C1 = '=""'
